$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"

# --- Simple numeric value updates ---
$ws.Range("L15").Value = 125
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 41
$ws.Range("K16").Value = -19.607843137254
$ws.Range("L16").Value = 17.142857142857
$ws.Range("M16").Value = -32.786885245901
$ws.Range("N16").Value = -83.534136546184
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 140
$ws.Range("K17").Value = -26.428571428571
$ws.Range("L17").Value = -9.649122807017
$ws.Range("M17").Value = 32.051282051282
$ws.Range("N17").Value = -49.261083743842
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 54
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = 5.882352941176
$ws.Range("L18").Value = 38.461538461538
$ws.Range("M18").Value = 1.886792452830
$ws.Range("N18").Value = -85.900783289817
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = -42.857142857142
$ws.Range("I19").Value = 118
$ws.Range("J19").Value = 147
$ws.Range("K19").Value = -19.727891156462
$ws.Range("L19").Value = -32.954545454545
$ws.Range("M19").Value = 11.320754716981
$ws.Range("N19").Value = -34.806629834254
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = -48.076923076923
$ws.Range("L20").Value = -28.947368421052
$ws.Range("N20").Value = -90.784982935153
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = -44.444444444444
$ws.Range("F21").Value = 28
$ws.Range("G21").Value = 33
$ws.Range("H21").Value = -15.151515151515
$ws.Range("I21").Value = 354
$ws.Range("J21").Value = 450
$ws.Range("K21").Value = -21.333333333333
$ws.Range("L21").Value = -12.807881773399
$ws.Range("M21").Value = 8.256880733944
$ws.Range("N21").Value = -73.303167420814
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -60
$ws.Range("I23").Value = 26
$ws.Range("K23").Value = 23.809523809523
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = 36.842105263157
$ws.Range("C24").Value = 8
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 31
$ws.Range("H24").Value = 22.580645161290
$ws.Range("I24").Value = 415
$ws.Range("J24").Value = 457
$ws.Range("K24").Value = -9.190371991247
$ws.Range("L24").Value = -3.037383177570
$ws.Range("M24").Value = 50.362318840579
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 204
$ws.Range("J25").Value = 234
$ws.Range("K25").Value = -12.820512820512
$ws.Range("L25").Value = -7.692307692307
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 206
$ws.Range("J26").Value = 237
$ws.Range("K26").Value = -13.080168776371
$ws.Range("L26").Value = 3
$ws.Range("M26").Value = -33.974358974359
$ws.Range("L27").Value = 57.142857142857

# --- Text -> Number conversions (restore numeric style) ---
$ws.Range("C16").Value = 2
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = -50
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Number -> Text conversions (use shared text style like A16) ---
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "0"
$ws.Range("A16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("A16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("A16").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("A16").Copy()
$ws.Range("H28").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Edit complete"